$d = $word.ActiveDocument

# --- Step 1: split/extend the final paragraph (Sprint 1 Retrospective ...) ---
# Replace its content with itself (plus the new pPr/rPr rStyle) followed by
# six brand-new paragraphs: Live URL heading, URL-hyperlink placeholder,
# Login Credentials heading, and the two credential paragraphs. InsertXML
# keeps the original paragraph-mark as a trailing empty paragraph, which
# becomes the new final (bookmark) paragraph.
$p = $d.Paragraphs.Last
$rng = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr></w:pPr><w:r><w:t>Sprint 1 Retrospective</w:t></w:r><w:r><w:t xml:space="preserve"> : </w:t></w:r><w:hyperlink r:id="rId8" w:tgtFrame="_blank" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://youtu.be/kxd9nJxSA80</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>Live URL</w:t></w:r></w:p><w:p><w:r><w:t>ZZZLIVEURLPLACEHOLDERZZZ</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Login Credentials</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Demand management group </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>user</w:t></w:r><w:r><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> userdm1/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>passwd</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Application group </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>user</w:t></w:r><w:r><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> userapp1/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>passwd</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# --- Step 2: InsertXML drops rStyle on *runs* (pPr/rPr survives), so re-apply
#     the Hyperlink character style to the run(s) that need it. ---
$f = $d.Content
$f.Find.Execute("https://youtu.be/kxd9nJxSA80", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$f.Style = "Hyperlink"

# --- Step 3: turn the placeholder paragraph into the real hyperlink run
#     pointing at the live URL (this also creates the external relationship). ---
$g = $d.Content
$g.Find.Execute("ZZZLIVEURLPLACEHOLDERZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($g, "http://52.5.93.209:8082/wpm/index.jsp", $null, $null, "http://52.5.93.209:8082/wpm/index.jsp") | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
